$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Column D (ID) updates
$ws.Range("D6").Value = 60000305
$ws.Range("D7").Value = 60000305
$ws.Range("D8").Value = 60000305
$ws.Range("D9").Value = 60000305
$ws.Range("D10").Value = 60000305
$ws.Range("D16").Value = 60000305
$ws.Range("D17").Value = 60000305
$ws.Range("D20").Value = 60000306
$ws.Range("D22").Value = 60000307
$ws.Range("D24").Value = 60000308
$ws.Range("D26").Value = 60000309

# Column E (sequence) updates
$ws.Range("E11").Value = 235
$ws.Range("E12").Value = 235
$ws.Range("E13").Value = 235
$ws.Range("E14").Value = 235
$ws.Range("E15").Value = 235
$ws.Range("E18").Value = 235
$ws.Range("E19").Value = 235
$ws.Range("E21").Value = 236
$ws.Range("E23").Value = 237
$ws.Range("E25").Value = 238
$ws.Range("E27").Value = 239
